$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the executed-run timestamps in column B (rows 2-5) to the
# latest RAD test-run values, matching the new shared-string entries
# appended for this commit ("Added RAD Test Cases and data for MRF.").
$ws.Range("B2").Value = "Mon Oct 02 16:44:08 EDT 2023"
$ws.Range("B3").Value = "Mon Oct 02 16:44:21 EDT 2023"
$ws.Range("B4").Value = "Mon Oct 02 16:44:34 EDT 2023"
$ws.Range("B5").Value = "Mon Oct 02 16:44:46 EDT 2023"
